$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000000091429589899050256462942
$ws.Range("C2").Value = 0.002658071450198252073243709859
$ws.Range("D2").Value = 0.721094517987026506489200983196
$ws.Range("E2").Value = 13.863846470800680776847002562135
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 14.58759915166750076309654104989
